# Auto-generated edit script applying F-column ("想去人数") updates
# as described by the diff, per sheet: 展览 (sheet1), 本地生活 (sheet3), 全部类型 (sheet4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7986
$ws.Range("F7").Value = 65
$ws.Range("F8").Value = 6897
$ws.Range("F9").Value = 6897
$ws.Range("F11").Value = 508
$ws.Range("F12").Value = 476
$ws.Range("F14").Value = 687
$ws.Range("F21").Value = 11265
$ws.Range("F22").Value = 90
$ws.Range("F23").Value = 2144
$ws.Range("F25").Value = 2937
$ws.Range("F28").Value = 2559
$ws.Range("F29").Value = 97
$ws.Range("F31").Value = 261
$ws.Range("F34").Value = 2307
$ws.Range("F36").Value = 1567
$ws.Range("F38").Value = 75
$ws.Range("F39").Value = 5687
$ws.Range("F40").Value = 1748
$ws.Range("F41").Value = 1239
$ws.Range("F42").Value = 813
$ws.Range("F43").Value = 152
$ws.Range("F46").Value = 1099
$ws.Range("F47").Value = 1055
$ws.Range("F48").Value = 1486

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 171
$ws.Range("F3").Value = 292

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 171
$ws.Range("F5").Value = 292
$ws.Range("F8").Value = 7986
$ws.Range("F11").Value = 65
$ws.Range("F12").Value = 6897
$ws.Range("F14").Value = 508
$ws.Range("F15").Value = 476
$ws.Range("F16").Value = 687
$ws.Range("F24").Value = 11265
$ws.Range("F25").Value = 90
$ws.Range("F26").Value = 2144
$ws.Range("F27").Value = 2937
$ws.Range("F28").Value = 2560
$ws.Range("F30").Value = 261
$ws.Range("F33").Value = 2307
$ws.Range("F35").Value = 1567
$ws.Range("F37").Value = 75
$ws.Range("F38").Value = 5687
$ws.Range("F40").Value = 1748
$ws.Range("F42").Value = 1239
$ws.Range("F43").Value = 813
$ws.Range("F44").Value = 152
$ws.Range("F46").Value = 1099
$ws.Range("F47").Value = 1055
$ws.Range("F48").Value = 1486
